$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for the "380 kV" case (rows 2-25, columns
# B, C, D, E, F, H, I, L, M), per the source XLSX diff.

$ws.Range("B2").Value = 18.45799994421992
$ws.Range("C2").Value = 9.160868140124864
$ws.Range("D2").Value = 7.74780320533904
$ws.Range("E2").Value = 9.868800480874159
$ws.Range("F2").Value = 40.75843688601577
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 32.27767050510458
$ws.Range("L2").Value = 10.35085250042845
$ws.Range("M2").Value = 16.68056971028239
$ws.Range("B3").Value = 18.12031132147241
$ws.Range("C3").Value = 8.549156298291283
$ws.Range("D3").Value = 7.761444272748255
$ws.Range("E3").Value = 9.853265667986747
$ws.Range("F3").Value = 40.40164030969925
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 32.17800951951012
$ws.Range("L3").Value = 10.36146805487871
$ws.Range("M3").Value = 16.63224542646556
$ws.Range("B4").Value = 17.91599688949036
$ws.Range("C4").Value = 8.180485582758831
$ws.Range("D4").Value = 7.770545946825372
$ws.Range("E4").Value = 9.843458967855673
$ws.Range("F4").Value = 40.19174032264959
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 32.12330283803291
$ws.Range("L4").Value = 10.36956841078441
$ws.Range("M4").Value = 16.60650184841728
$ws.Range("B5").Value = 17.83362574962241
$ws.Range("C5").Value = 8.032866124750848
$ws.Range("D5").Value = 7.77443683348076
$ws.Range("E5").Value = 9.839394905585182
$ws.Range("F5").Value = 40.10858002655517
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 32.1026454165824
$ws.Range("L5").Value = 10.37326714411945
$ws.Range("M5").Value = 16.59700625996366
$ws.Range("B6").Value = 17.82000546900203
$ws.Range("C6").Value = 8.008117276074433
$ws.Range("D6").Value = 7.77509387756264
$ws.Range("E6").Value = 9.838715970453491
$ws.Range("F6").Value = 40.09491670895471
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 32.0993141921601
$ws.Range("L6").Value = 10.37390533536549
$ws.Range("M6").Value = 16.59548982568691
$ws.Range("B7").Value = 17.91488223754037
$ws.Range("C7").Value = 8.178510700732799
$ws.Range("D7").Value = 7.770597685174736
$ws.Range("E7").Value = 9.843404433179623
$ws.Range("F7").Value = 40.19060908834183
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 32.12301761548833
$ws.Range("L7").Value = 10.36961668285732
$ws.Range("M7").Value = 16.60636974932464
$ws.Range("B8").Value = 18.34101288855101
$ws.Range("C8").Value = 8.954773425292279
$ws.Range("D8").Value = 7.752355554885988
$ws.Range("E8").Value = 9.863499574267383
$ws.Range("F8").Value = 40.63355198149902
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 32.241964933774
$ws.Range("L8").Value = 10.35418426784245
$ws.Range("M8").Value = 16.6630970896989
$ws.Range("B9").Value = 19.1951050967572
$ws.Range("C9").Value = 10.35285892843316
$ws.Range("D9").Value = 7.722373731891127
$ws.Range("E9").Value = 9.900791979787225
$ws.Range("F9").Value = 41.57170681545403
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 32.52637909115693
$ws.Range("L9").Value = 10.33648063554044
$ws.Range("M9").Value = 16.80512984684274
$ws.Range("B10").Value = 19.82651826865609
$ws.Range("C10").Value = 11.26912385133709
$ws.Range("D10").Value = 7.703918228107545
$ws.Range("E10").Value = 9.926923476377512
$ws.Range("F10").Value = 42.29871886166525
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 32.76601570514836
$ws.Range("L10").Value = 10.3311331811732
$ws.Range("M10").Value = 16.92768607581563
$ws.Range("B11").Value = 20.1131241821836
$ws.Range("C11").Value = 11.66213646970422
$ws.Range("D11").Value = 7.696307075902002
$ws.Range("E11").Value = 9.938540951006672
$ws.Range("F11").Value = 42.63659173468183
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 32.88154896831794
$ws.Range("L11").Value = 10.33036276751773
$ws.Range("M11").Value = 16.98725223177431
$ws.Range("B12").Value = 20.22144635513911
$ws.Range("C12").Value = 11.80755800301248
$ws.Range("D12").Value = 7.69353848452628
$ws.Range("E12").Value = 9.942901724526745
$ws.Range("F12").Value = 42.76547086945218
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 32.92622100209537
$ws.Range("L12").Value = 10.33030979743408
$ws.Range("M12").Value = 17.01034358829608
$ws.Range("B13").Value = 20.19812838030639
$ws.Range("C13").Value = 11.77638996271146
$ws.Range("D13").Value = 7.694129685638742
$ws.Range("E13").Value = 9.941964262029819
$ws.Range("F13").Value = 42.73767442036534
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 32.9165593365738
$ws.Range("L13").Value = 10.33031059108105
$ws.Range("M13").Value = 17.00534688521664
$ws.Range("B14").Value = 20.12204076865901
$ws.Range("C14").Value = 11.67416849751456
$ws.Range("D14").Value = 7.696077021184403
$ws.Range("E14").Value = 9.938900484602392
$ws.Range("F14").Value = 42.64717652579163
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 32.88520577879081
$ws.Range("L14").Value = 10.33035362609293
$ws.Range("M14").Value = 16.98914133138899
$ws.Range("B15").Value = 20.07540421824601
$ws.Range("C15").Value = 11.61111217471687
$ws.Range("D15").Value = 7.697284636567667
$ws.Range("E15").Value = 9.937018825885913
$ws.Range("F15").Value = 42.59186276016681
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 32.86612041048955
$ws.Range("L15").Value = 10.33041107214066
$ws.Range("M15").Value = 16.97928420926308
$ws.Range("B16").Value = 19.80776497278789
$ws.Range("C16").Value = 11.24296222111865
$ws.Range("D16").Value = 7.704431486034617
$ws.Range("E16").Value = 9.926158841743034
$ws.Range("F16").Value = 42.27677395295873
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 32.75859525900321
$ws.Range("L16").Value = 10.33121695732677
$ws.Range("M16").Value = 16.92386890576596
$ws.Range("B17").Value = 19.64332988066567
$ws.Range("C17").Value = 11.0110315329405
$ws.Range("D17").Value = 7.709017303638934
$ws.Range("E17").Value = 9.919427681537735
$ws.Range("F17").Value = 42.08524226062007
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 32.69429171160598
$ws.Range("L17").Value = 10.33213691385737
$ws.Range("M17").Value = 16.89084104558545
$ws.Range("B18").Value = 19.54869893248694
$ws.Range("C18").Value = 10.87538977111624
$ws.Range("D18").Value = 7.711728723245162
$ws.Range("E18").Value = 9.915530676549551
$ws.Range("F18").Value = 41.97575720334574
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 32.65792087817659
$ws.Range("L18").Value = 10.33282252032158
$ws.Range("M18").Value = 16.8722042235064
$ws.Range("B19").Value = 19.51665326715879
$ws.Range("C19").Value = 10.82907827235406
$ws.Range("D19").Value = 7.71265941062839
$ws.Range("E19").Value = 9.914206842814597
$ws.Range("F19").Value = 41.93880689339175
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 32.64571239329649
$ws.Range("L19").Value = 10.33308153521737
$ws.Range("M19").Value = 16.86595633441652
$ws.Range("B20").Value = 19.66084056312134
$ws.Range("C20").Value = 11.03595278723
$ws.Range("D20").Value = 7.708521493957093
$ws.Range("E20").Value = 9.920146853736989
$ws.Range("F20").Value = 42.10556155024061
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 32.70107340602839
$ws.Range("L20").Value = 10.33202279020861
$ws.Range("M20").Value = 16.89431976219256
$ws.Range("B21").Value = 20.14439614516355
$ws.Range("C21").Value = 11.7042856497223
$ws.Range("D21").Value = 7.695501952427811
$ws.Range("E21").Value = 9.939801433137555
$ws.Range("F21").Value = 42.67373336214639
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 32.89439018624474
$ws.Range("L21").Value = 10.33033450796935
$ws.Range("M21").Value = 16.99388688970299
$ws.Range("B22").Value = 20.45915674971877
$ws.Range("C22").Value = 12.12124836590296
$ws.Range("D22").Value = 7.687655374302667
$ws.Range("E22").Value = 9.952422498687669
$ws.Range("F22").Value = 43.05046587284417
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 33.02610015720075
$ws.Range("L22").Value = 10.33062266783917
$ws.Range("M22").Value = 17.06207094498924
$ws.Range("B23").Value = 20.29131724473485
$ws.Range("C23").Value = 11.90051594257458
$ws.Range("D23").Value = 7.69178236179761
$ws.Range("E23").Value = 9.945706827119064
$ws.Range("F23").Value = 42.8489341464383
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 32.95531864238402
$ws.Range("L23").Value = 10.33034164957241
$ws.Range("M23").Value = 17.02539986773169
$ws.Range("B24").Value = 19.65292427046244
$ws.Range("C24").Value = 11.02469305892804
$ws.Range("D24").Value = 7.7087454160534
$ws.Range("E24").Value = 9.919821800640525
$ws.Range("F24").Value = 42.09637323287201
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 32.69800553776122
$ws.Range("L24").Value = 10.33207389737184
$ws.Range("M24").Value = 16.89274593906374
$ws.Range("B25").Value = 18.96284380967276
$ws.Range("C25").Value = 9.994205526449971
$ws.Range("D25").Value = 7.729860346166999
$ws.Range("E25").Value = 9.890927500743627
$ws.Range("F25").Value = 41.31093064449767
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 32.44400619024844
$ws.Range("L25").Value = 10.33992475564372
$ws.Range("M25").Value = 16.76346434626941
